$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns E..T (F and L stay at 1, unchanged)
$data = @{
    2  = @{ E=3; G=1.542357666666667; H=4.627073; I=0.1736642262104688; J=0.1736642262104688; K=3; M=405.24646; N=1215.73938; O=0.971171031955694; P=0.9711710319556939; Q=625.0349844705268; R=5625.31486023474; S=0.168657665782608; T=0.168657665782608 }
    3  = @{ E=3; G=1.542357666666667; H=4.627073; I=0.1736642262104688; J=0.1736642262104688; K=3; M=5.631177666666667; N=16.893533; O=0.01349508796612936; P=0.01349508796612936; Q=8.685290046545445; R=78.16761041890901; S=0.002343614009280065; T=0.002343614009280065 }
    4  = @{ E=3; G=1.542357666666667; H=4.627073; I=0.1736642262104688; J=0.1736642262104688; K=3; M=1.958728333333333; N=5.876185; O=0.004694082255041018; P=0.004694082255041017; Q=3.021059661833889; R=27.189536956505; S=0.0008151941625899908; T=0.0008151941625899907 }
    5  = @{ E=3; G=1.542357666666667; H=4.627073; I=0.1736642262104688; J=0.1736642262104688; K=3; M=3.263573; N=9.790718999999999; O=0.007821135706583939; P=0.007821135706583937; Q=5.033596837276333; R=45.302371535487; S=0.001358251480570968; T=0.001358251480570968 }
    6  = @{ E=3; G=1.542357666666667; H=4.627073; I=0.1736642262104688; J=0.1736642262104688; K=3; M=1.176160333333333; N=3.528481; O=0.002818662116551706; P=0.002818662116551706; Q=1.814059907345889; R=16.326539166113; S=0.0004895007754197143; T=0.0004895007754197141 }
    7  = @{ E=3; G=7.338905333333334; H=22.016716; I=0.8263357737895313; J=0.8263357737895313; K=3; M=405.24646; N=1215.73938; O=0.971171031955694; P=0.9711710319556939; Q=2974.065406608453; R=26766.58865947608; S=0.802513366173086; T=0.8025133661730859 }
    8  = @{ E=3; G=7.338905333333334; H=22.016716; I=0.8263357737895313; J=0.8263357737895313; K=3; M=5.631177666666667; N=16.893533; O=0.01349508796612936; P=0.01349508796612936; Q=41.32667981084756; R=371.9401182976281; S=0.0111514739568493; T=0.0111514739568493 }
    9  = @{ E=3; G=7.338905333333334; H=22.016716; I=0.8263357737895313; J=0.8263357737895313; K=3; M=1.958728333333333; N=5.876185; O=0.004694082255041018; P=0.004694082255041017; Q=14.37492181205111; R=129.37429630846; S=0.003878888092451028; T=0.003878888092451027 }
    10 = @{ E=3; G=7.338905333333334; H=22.016716; I=0.8263357737895313; J=0.8263357737895313; K=3; M=3.263573; N=9.790718999999999; O=0.007821135706583939; P=0.007821135706583937; Q=23.95105329542266; R=215.559479658804; S=0.006462884226012972; T=0.00646288422601297 }
    11 = @{ E=3; G=7.338905333333334; H=22.016716; I=0.8263357737895313; J=0.8263357737895313; K=3; M=1.176160333333333; N=3.528481; O=0.002818662116551706; P=0.002818662116551706; Q=8.631729343155111; R=77.68556408839601; S=0.002329161341131992; T=0.002329161341131992 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
